# Update Ngf-Sort1 NATMI TPM output with recomputed values (new TPM run).
# Sheet has a single data table (rows 2-10, cols A-T); only numeric result
# columns (E-T) change per the updated TPM calculation - identifying
# columns A-D (clusters/ligand/receptor) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.322531
$ws.Range("H2").Value = 0.967593
$ws.Range("I2").Value = 0.01892149513432853
$ws.Range("J2").Value = 0.01892149513432853
$ws.Range("M2").Value = 0.9386610000000001
$ws.Range("N2").Value = 2.815983
$ws.Range("O2").Value = 0.04181245246793033
$ws.Range("P2").Value = 0.04181245246793032
$ws.Range("Q2").Value = 0.302747270991
$ws.Range("R2").Value = 2.724725438919
$ws.Range("S2").Value = 0.0007911541159262868
$ws.Range("T2").Value = 0.0007911541159262865
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.322531
$ws.Range("H3").Value = 0.967593
$ws.Range("I3").Value = 0.01892149513432853
$ws.Range("J3").Value = 0.01892149513432853
$ws.Range("O3").Value = 0.1106393125456779
$ws.Range("P3").Value = 0.1106393125456779
$ws.Range("Q3").Value = 0.801095079587
$ws.Range("R3").Value = 7.209855716283
$ws.Range("S3").Value = 0.002093461213998499
$ws.Range("T3").Value = 0.002093461213998498
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.322531
$ws.Range("H4").Value = 0.967593
$ws.Range("I4").Value = 0.01892149513432853
$ws.Range("J4").Value = 0.01892149513432853
$ws.Range("O4").Value = 0.8475482349863918
$ws.Range("P4").Value = 0.8475482349863918
$ws.Range("Q4").Value = 6.136758310748999
$ws.Range("R4").Value = 55.230824796741
$ws.Range("S4").Value = 0.01603687980440375
$ws.Range("T4").Value = 0.01603687980440375
$ws.Range("I5").Value = 0.0261208867009986
$ws.Range("J5").Value = 0.0261208867009986
$ws.Range("M5").Value = 0.9386610000000001
$ws.Range("N5").Value = 2.815983
$ws.Range("O5").Value = 0.04181245246793033
$ws.Range("P5").Value = 0.04181245246793032
$ws.Range("Q5").Value = 0.41793881025
$ws.Range("R5").Value = 3.76144929225
$ws.Range("S5").Value = 0.001092178333605697
$ws.Range("T5").Value = 0.001092178333605697
$ws.Range("I6").Value = 0.0261208867009986
$ws.Range("J6").Value = 0.0261208867009986
$ws.Range("O6").Value = 0.1106393125456779
$ws.Range("P6").Value = 0.1106393125456779
$ws.Range("S6").Value = 0.002889996947682025
$ws.Range("T6").Value = 0.002889996947682025
$ws.Range("I7").Value = 0.0261208867009986
$ws.Range("J7").Value = 0.0261208867009986
$ws.Range("O7").Value = 0.8475482349863918
$ws.Range("P7").Value = 0.8475482349863918
$ws.Range("S7").Value = 0.02213871141971088
$ws.Range("T7").Value = 0.02213871141971087
$ws.Range("I8").Value = 0.954957618164673
$ws.Range("J8").Value = 0.954957618164673
$ws.Range("M8").Value = 0.9386610000000001
$ws.Range("N8").Value = 2.815983
$ws.Range("O8").Value = 0.04181245246793033
$ws.Range("P8").Value = 0.04181245246793032
$ws.Range("Q8").Value = 15.279490904865
$ws.Range("R8").Value = 137.515418143785
$ws.Range("S8").Value = 0.03992912001839834
$ws.Range("T8").Value = 0.03992912001839834
$ws.Range("I9").Value = 0.954957618164673
$ws.Range("J9").Value = 0.954957618164673
$ws.Range("O9").Value = 0.1106393125456779
$ws.Range("P9").Value = 0.1106393125456779
$ws.Range("S9").Value = 0.1056558543839974
$ws.Range("T9").Value = 0.1056558543839974
$ws.Range("I10").Value = 0.954957618164673
$ws.Range("J10").Value = 0.954957618164673
$ws.Range("O10").Value = 0.8475482349863918
$ws.Range("P10").Value = 0.8475482349863918
$ws.Range("S10").Value = 0.8093726437622772
$ws.Range("T10").Value = 0.8093726437622772
